$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 843.375
$ws.Range("J17").Value = 999
$ws.Range("L17").Value = 2997
$ws.Range("N17").Value = -3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H16").Value = 5795.5835
$ws.Range("I16").Value = 2859.5
$ws.Range("J16").Value = 8731.666999999999
$ws.Range("K16").Value = 2859.5
$ws.Range("L16").Value = 8731.666999999999
$ws.Range("M16").Value = -2572.5
$ws.Range("N16").Value = -9305.666999999999
$ws.Range("H19").Value = 59001.75
$ws.Range("I19").Value = 59001.75
$ws.Range("K19").Value = 59001.75
$ws.Range("M19").Value = -58772.75
$ws.Range("H46").Value = 9699.4
$ws.Range("I46").Value = 9998
$ws.Range("J46").Value = 9624.75
$ws.Range("K46").Value = 9998
$ws.Range("L46").Value = 9624.75
$ws.Range("M46").Value = -9679
$ws.Range("N46").Value = -10262.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1914.7812
$ws.Range("I20").Value = 1961.1923
$ws.Range("J20").Value = 1713.6666
$ws.Range("K20").Value = 1961.1923
$ws.Range("L20").Value = 1713.6666
$ws.Range("M20").Value = -1714.1923
$ws.Range("N20").Value = -2207.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 4922.5
$ws.Range("I12").Value = 3464.2
$ws.Range("J12").Value = 7353
$ws.Range("K12").Value = 3464.2
$ws.Range("L12").Value = 7353
$ws.Range("M12").Value = -3294.2
$ws.Range("N12").Value = -7693
$ws.Range("H25").Value = 7898.6665
$ws.Range("I25").Value = 6340
$ws.Range("J25").Value = 9847
$ws.Range("K25").Value = 6340
$ws.Range("L25").Value = 9847
$ws.Range("M25").Value = -6166
$ws.Range("N25").Value = -10195
$ws.Range("H58").Value = 2770.025
$ws.Range("I58").Value = 2412
$ws.Range("J58").Value = 4457.857
$ws.Range("K58").Value = 2412
$ws.Range("L58").Value = 4457.857
$ws.Range("M58").Value = -2209
$ws.Range("N58").Value = -4863.857
$ws.Range("H136").Value = 2770.025
$ws.Range("I136").Value = 2412
$ws.Range("J136").Value = 4457.857
$ws.Range("K136").Value = 7236
$ws.Range("L136").Value = 13373.571
$ws.Range("M136").Value = -4686
$ws.Range("N136").Value = -18473.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 2034.1818
$ws.Range("I13").Value = 930.1667
$ws.Range("K13").Value = 2790.5001
$ws.Range("M13").Value = -2622.5001
$ws.Range("H137").Value = 3079.7273
$ws.Range("I137").Value = 1375.5454
$ws.Range("K137").Value = 4126.6362
$ws.Range("M137").Value = 973.3638000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 4000000
$ws.Range("I10").Value = 4000000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 4000000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -3999831
$ws.Range("N10").ClearContents()
$ws.Range("H11").Value = 5800400
$ws.Range("J11").Value = 3000666.8
$ws.Range("L11").Value = 3000666.8
$ws.Range("N11").Value = -3000944.8
$ws.Range("H12").Value = 2509999.5
$ws.Range("I12").Value = 2509999.5
$ws.Range("K12").Value = 2509999.5
$ws.Range("M12").Value = -2509859.5
$ws.Range("H18").Value = 8794.714
$ws.Range("I18").Value = 8794.714
$ws.Range("K18").Value = 8794.714
$ws.Range("M18").Value = -8501.714
$ws.Range("H36").Value = 44604476
$ws.Range("I36").Value = 70089384
$ws.Range("J36").Value = 5881.25
$ws.Range("K36").Value = 70089384
$ws.Range("L36").Value = 5881.25
$ws.Range("M36").Value = -70088899
$ws.Range("N36").Value = -6851.25
$ws.Range("H134").Value = 62442.5
$ws.Range("J134").Value = 62442.5
$ws.Range("L134").Value = 187327.5
$ws.Range("N134").Value = -192397.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 3499.7144
$ws.Range("I12").Value = 699.6667
$ws.Range("K12").Value = 699.6667
$ws.Range("M12").Value = -529.6667
$ws.Range("H23").Value = 2930.6667
$ws.Range("I23").Value = 2930.6667
$ws.Range("K23").Value = 2930.6667
$ws.Range("M23").Value = -2700.6667
$ws.Range("H58").Value = 11023.25
$ws.Range("I58").Value = 4996.5
$ws.Range("J58").Value = 17050
$ws.Range("K58").Value = 4996.5
$ws.Range("L58").Value = 17050
$ws.Range("M58").Value = -4736.5
$ws.Range("N58").Value = -17570

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 7499.5
$ws.Range("I6").Value = 7499.5
$ws.Range("K6").Value = 7499.5
$ws.Range("M6").Value = -7384.5
$ws.Range("H7").Value = 4999
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 4999
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 4999
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -5225
$ws.Range("H10").Value = 43999.4
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 43999.4
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 43999.4
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -44337.4
$ws.Range("H30").Value = 1999.5
$ws.Range("I30").Value = 1999.5
$ws.Range("K30").Value = 1999.5
$ws.Range("M30").Value = -1892.5
$ws.Range("H131").Value = 65999.39999999999
$ws.Range("J131").Value = 65999.39999999999
$ws.Range("L131").Value = 65999.39999999999
$ws.Range("N131").Value = -76079.39999999999
$ws.Range("H132").Value = 5458.0713
$ws.Range("I132").Value = 5569.846
$ws.Range("K132").Value = 16709.538
$ws.Range("M132").Value = -14179.538
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 3967.111
$ws.Range("I136").Value = 4374.4546
$ws.Range("J136").Value = 3023.7896
$ws.Range("K136").Value = 13123.3638
$ws.Range("L136").Value = 9071.3688
$ws.Range("M136").Value = -10573.3638
$ws.Range("N136").Value = -14171.3688
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 115000
$ws.Range("J139").Value = 115000
$ws.Range("L139").Value = 115000
$ws.Range("N139").Value = -125280
$ws.Range("H141").Value = 120749.75
$ws.Range("J141").Value = 120749.75
$ws.Range("L141").Value = 120749.75
$ws.Range("N141").Value = -131109.75
